{"js": "// Helper that locates a paragraph whose text contains the given marker\n// and replaces its whole range text with the new text.\nasync function replaceParagraphContaining(marker, newText) {\n  const body = context.document.body;\n  const paragraphs = body.paragraphs;\n  paragraphs.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < paragraphs.items.length; i++) {\n    const para = paragraphs.items[i];\n    para.load(\"text\");\n  }\n  await context.sync();\n\n  for (let i = 0; i < paragraphs.items.length; i++) {\n    const para = paragraphs.items[i];\n    if (para.text.indexOf(marker) !== -1) {\n      const range = para.getRange();\n      range.insertText(newText, Word.InsertLocation.replace);\n      return true;\n    }\n  }\n  return false;\n}\n\nawait replaceParagraphContaining(\n  \"\u03a6.350.2/1/32958\",\n  \"\u03a4\u03b7 \u03bc\u03b5 \u03b1\u03c1\u03b9\u03b8. \u03a6.351.1/11/48020/\u03953/28-3-2019 (\u0391\u0394\u0391: \u03a9\u03a9\u03a4\u03974653\u03a0\u03a3-\u0392\u03943) \u03a5\u03c0\u03bf\u03c5\u03c1\u03b3\u03b9\u03ba\u03ae \u0391\u03c0\u03cc\u03c6\u03b1\u03c3\u03b7 \u03bc\u03b5 \u03b8\u03ad\u03bc\u03b1: \u00ab\u03a4\u03bf\u03c0\u03bf\u03b8\u03ad\u03c4\u03b7\u03c3\u03b7 \u03a0\u03b5\u03c1\u03b9\u03c6\u03b5\u03c1\u03b5\u03b9\u03b1\u03ba\u03ce\u03bd \u0394\u03b9\u03b5\u03c5\u03b8\u03c5\u03bd\u03c4\u03ce\u03bd \u0395\u03ba\u03c0\u03b1\u03af\u03b4\u03b5\u03c5\u03c3\u03b7\u03c2\u00bb\"\n);\n\nawait replaceParagraphContaining(\n  \"\u03a6.353.1/324/105657\",\n  \"\u03a4\u03b7 \u03bc\u03b5 \u03b1\u03c1\u03b9\u03b8. \u03a6.353.1/324/105657/\u03941/8-10-2002 \u0391\u03c0\u03cc\u03c6\u03b1\u03c3\u03b7 \u03a5\u03c6\u03c5\u03c0\u03bf\u03c5\u03c1\u03b3\u03bf\u03cd \u0395\u03b8\u03bd\u03b9\u03ba\u03ae\u03c2 \u03a0\u03b1\u03b9\u03b4\u03b5\u03af\u03b1\u03c2 \u03ba\u03b1\u03b9 \u0398\u03c1\u03b7\u03c3\u03ba\u03b5\u03c5\u03bc\u03ac\u03c4\u03c9\u03bd \u00ab\u039a\u03b1\u03b8\u03bf\u03c1\u03b9\u03c3\u03bc\u03cc\u03c2 \u03c4\u03c9\u03bd \u03b5\u03b9\u03b4\u03b9\u03ba\u03cc\u03c4\u03b5\u03c1\u03c9\u03bd \u03ba\u03b1\u03b8\u03b7\u03ba\u03cc\u03bd\u03c4\u03c9\u03bd \u03ba\u03b1\u03b9 \u03b1\u03c1\u03bc\u03bf\u03b4\u03b9\u03bf\u03c4\u03ae\u03c4\u03c9\u03bd \u03c4\u03c9\u03bd \u03a0\u03c1\u03bf\u03ca\u03c3\u03c4\u03b1\u03bc\u03ad\u03bd\u03c9\u03bd \u03c4\u03c9\u03bd \u03a0\u03b5\u03c1\u03b9\u03c6\u03b5\u03c1\u03b5\u03b9\u03b1\u03ba\u03ce\u03bd \u03c5\u03c0\u03b7\u03c1\u03b5\u03c3\u03b9\u03ce\u03bd \u03a0\u03c1\u03c9\u03c4\u03bf\u03b2\u03ac\u03b8\u03bc\u03b9\u03b1\u03c2 \u03ba\u03b1\u03b9 \u0394\u03b5\u03c5\u03c4\u03b5\u03c1\u03bf\u03b2\u03ac\u03b8\u03bc\u03b9\u03b1\u03c2 \u03b5\u03ba\u03c0\u03b1\u03af\u03b4\u03b5\u03c5\u03c3\u03b7\u03c2, \u03c4\u03c9\u03bd \u0394\u03b9\u03b5\u03c5\u03b8\u03c5\u03bd\u03c4\u03ce\u03bd \u03ba\u03b1\u03b9 \u03c5\u03c0\u03bf\u03b4\u03b9\u03b5\u03c5\u03b8\u03c5\u03bd\u03c4\u03ce\u03bd \u03c4\u03c9\u03bd \u03c3\u03c7\u03bf\u03bb\u03b9\u03ba\u03ce\u03bd \u03bc\u03bf\u03bd\u03ac\u03b4\u03c9\u03bd \u03ba\u03b1\u03b9 \u03a3\u0395\u039a \u03ba\u03b1\u03b9 \u03c4\u03c9\u03bd \u03c3\u03c5\u03bb\u03bb\u03cc\u03b3\u03c9\u03bd \u03c4\u03c9\u03bd \u03b4\u03b9\u03b4\u03b1\u03c3\u03ba\u03cc\u03bd\u03c4\u03c9\u03bd\u00bb (\u03a6\u0395\u039a 1340/\u03c4.\u0392\u0384/16-10-2002), \u03cc\u03c0\u03c9\u03c2 \u03c3\u03c5\u03bc\u03c0\u03bb\u03b7\u03c1\u03ce\u03b8\u03b7\u03ba\u03b5, \u03c4\u03c1\u03bf\u03c0\u03bf\u03c0\u03bf\u03b9\u03ae\u03b8\u03b7\u03ba\u03b5 \u03ba\u03b1\u03b9 \u03b9\u03c3\u03c7\u03cd\u03b5\u03b9 \u03c3\u03ae\u03bc\u03b5\u03c1\u03b1\"\n);\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Replace-ParagraphContaining($doc, $marker, $newText) {\n    foreach ($p in $doc.Paragraphs) {\n        $t = $p.Range.Text\n        if ($t -like \"*$marker*\") {\n            $r = $p.Range\n            $r.MoveEnd(1, -1) | Out-Null\n            $r.Text = $newText\n            return $true\n        }\n    }\n    return $false\n}\n\nReplace-ParagraphContaining $d \"350.2/1/32958\" \"\u03a4\u03b7 \u03bc\u03b5 \u03b1\u03c1\u03b9\u03b8. \u03a6.351.1/11/48020/\u03953/28-3-2019 (\u0391\u0394\u0391: \u03a9\u03a9\u03a4\u03974653\u03a0\u03a3-\u0392\u03943) \u03a5\u03c0\u03bf\u03c5\u03c1\u03b3\u03b9\u03ba\u03ae \u0391\u03c0\u03cc\u03c6\u03b1\u03c3\u03b7 \u03bc\u03b5 \u03b8\u03ad\u03bc\u03b1: \u00ab\u03a4\u03bf\u03c0\u03bf\u03b8\u03ad\u03c4\u03b7\u03c3\u03b7 \u03a0\u03b5\u03c1\u03b9\u03c6\u03b5\u03c1\u03b5\u03b9\u03b1\u03ba\u03ce\u03bd \u0394\u03b9\u03b5\u03c5\u03b8\u03c5\u03bd\u03c4\u03ce\u03bd \u0395\u03ba\u03c0\u03b1\u03af\u03b4\u03b5\u03c5\u03c3\u03b7\u03c2\u00bb\" | Out-Null\n\nReplace-ParagraphContaining $d \"353.1/324/105657\" \"\u03a4\u03b7 \u03bc\u03b5 \u03b1\u03c1\u03b9\u03b8. \u03a6.353.1/324/105657/\u03941/8-10-2002 \u0391\u03c0\u03cc\u03c6\u03b1\u03c3\u03b7 \u03a5\u03c6\u03c5\u03c0\u03bf\u03c5\u03c1\u03b3\u03bf\u03cd \u0395\u03b8\u03bd\u03b9\u03ba\u03ae\u03c2 \u03a0\u03b1\u03b9\u03b4\u03b5\u03af\u03b1\u03c2 \u03ba\u03b1\u03b9 \u0398\u03c1\u03b7\u03c3\u03ba\u03b5\u03c5\u03bc\u03ac\u03c4\u03c9\u03bd \u00ab\u039a\u03b1\u03b8\u03bf\u03c1\u03b9\u03c3\u03bc\u03cc\u03c2 \u03c4\u03c9\u03bd \u03b5\u03b9\u03b4\u03b9\u03ba\u03cc\u03c4\u03b5\u03c1\u03c9\u03bd \u03ba\u03b1\u03b8\u03b7\u03ba\u03cc\u03bd\u03c4\u03c9\u03bd \u03ba\u03b1\u03b9 \u03b1\u03c1\u03bc\u03bf\u03b4\u03b9\u03bf\u03c4\u03ae\u03c4\u03c9\u03bd \u03c4\u03c9\u03bd \u03a0\u03c1\u03bf\u03ca\u03c3\u03c4\u03b1\u03bc\u03ad\u03bd\u03c9\u03bd \u03c4\u03c9\u03bd \u03a0\u03b5\u03c1\u03b9\u03c6\u03b5\u03c1\u03b5\u03b9\u03b1\u03ba\u03ce\u03bd \u03c5\u03c0\u03b7\u03c1\u03b5\u03c3\u03b9\u03ce\u03bd \u03a0\u03c1\u03c9\u03c4\u03bf\u03b2\u03ac\u03b8\u03bc\u03b9\u03b1\u03c2 \u03ba\u03b1\u03b9 \u0394\u03b5\u03c5\u03c4\u03b5\u03c1\u03bf\u03b2\u03ac\u03b8\u03bc\u03b9\u03b1\u03c2 \u03b5\u03ba\u03c0\u03b1\u03af\u03b4\u03b5\u03c5\u03c3\u03b7\u03c2, \u03c4\u03c9\u03bd \u0394\u03b9\u03b5\u03c5\u03b8\u03c5\u03bd\u03c4\u03ce\u03bd \u03ba\u03b1\u03b9 \u03c5\u03c0\u03bf\u03b4\u03b9\u03b5\u03c5\u03b8\u03c5\u03bd\u03c4\u03ce\u03bd \u03c4\u03c9\u03bd \u03c3\u03c7\u03bf\u03bb\u03b9\u03ba\u03ce\u03bd \u03bc\u03bf\u03bd\u03ac\u03b4\u03c9\u03bd \u03ba\u03b1\u03b9 \u03a3\u0395\u039a \u03ba\u03b1\u03b9 \u03c4\u03c9\u03bd \u03c3\u03c5\u03bb\u03bb\u03cc\u03b3\u03c9\u03bd \u03c4\u03c9\u03bd \u03b4\u03b9\u03b4\u03b1\u03c3\u03ba\u03cc\u03bd\u03c4\u03c9\u03bd\u00bb (\u03a6\u0395\u039a 1340/\u03c4.\u0392\u0384/16-10-2002), \u03cc\u03c0\u03c9\u03c2 \u03c3\u03c5\u03bc\u03c0\u03bb\u03b7\u03c1\u03ce\u03b8\u03b7\u03ba\u03b5, \u03c4\u03c1\u03bf\u03c0\u03bf\u03c0\u03bf\u03b9\u03ae\u03b8\u03b7\u03ba\u03b5 \u03ba\u03b1\u03b9 \u03b9\u03c3\u03c7\u03cd\u03b5\u03b9 \u03c3\u03ae\u03bc\u03b5\u03c1\u03b1\" | Out-Null\n"}
